# Applies odds/value updates to Sheet1 as described in the commit diff.
# Each statement sets a single cell to its new numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.05
$ws.Range("K2").Value = 11
$ws.Range("N2").Value = 1.91
$ws.Range("O2").Value = 1.99
# Row 3
$ws.Range("AD3").Value = 201
$ws.Range("I3").Value = 6.25
$ws.Range("K3").Value = 17
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 2.1
$ws.Range("Y3").Value = 23
# Row 8
$ws.Range("AA8").Value = 8
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 101
$ws.Range("AE8").Value = 5
$ws.Range("AF8").Value = 6
$ws.Range("AG8").Value = 9
$ws.Range("AH8").Value = 10
$ws.Range("AI8").Value = 15
$ws.Range("AJ8").Value = 41
$ws.Range("G8").Value = 6.25
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 1.53
$ws.Range("J8").Value = 1.08
$ws.Range("K8").Value = 7.5
$ws.Range("L8").Value = 1.4
$ws.Range("M8").Value = 2.75
$ws.Range("N8").Value = 2.3
$ws.Range("O8").Value = 1.6
$ws.Range("P8").Value = 1.5
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 2.5
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 12
$ws.Range("U8").Value = 29
$ws.Range("W8").Value = 81
$ws.Range("Y8").Value = 67
$ws.Range("Z8").Value = 7.5
# Row 9
$ws.Range("AB9").Value = 15
$ws.Range("AD9").Value = 251
$ws.Range("AG9").Value = 8.5
$ws.Range("G9").Value = 4.75
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 1.83
$ws.Range("X9").Value = 41
# Row 11
$ws.Range("AG11").Value = 18.5
$ws.Range("R11").Value = 2.2
# Row 12
$ws.Range("AB12").Value = 18
$ws.Range("AE12").Value = 7.2
$ws.Range("AF12").Value = 14.5
$ws.Range("AG12").Value = 12
$ws.Range("AJ12").Value = 55
$ws.Range("M12").Value = 2.32
$ws.Range("N12").Value = 2.35
$ws.Range("R12").Value = 2
$ws.Range("T12").Value = 6.1
$ws.Range("U12").Value = 10
$ws.Range("V12").Value = 9.5
$ws.Range("X12").Value = 22
$ws.Range("Z12").Value = 6.6
# Row 14
$ws.Range("G14").Value = 2.55
$ws.Range("I14").Value = 2.9
$ws.Range("K14").Value = 8
$ws.Range("L14").Value = 1.44
$ws.Range("M14").Value = 2.63
$ws.Range("W14").Value = 26
# Row 17
$ws.Range("AA17").Value = 5.9
$ws.Range("AB17").Value = 13.5
$ws.Range("AC17").Value = 60
$ws.Range("AD17").Value = 450
$ws.Range("AE17").Value = 10.25
$ws.Range("AF17").Value = 22
$ws.Range("AG17").Value = 13
$ws.Range("AH17").Value = 65
$ws.Range("AI17").Value = 40
$ws.Range("AJ17").Value = 40
$ws.Range("G17").Value = 1.65
$ws.Range("H17").Value = 3.4
$ws.Range("I17").Value = 4.85
$ws.Range("N17").Value = 1.9
$ws.Range("O17").Value = 1.72
$ws.Range("T17").Value = 5.5
$ws.Range("U17").Value = 6.4
$ws.Range("V17").Value = 6.8
$ws.Range("W17").Value = 10.25
$ws.Range("X17").Value = 11.25
# Row 18
$ws.Range("AA18").Value = 7.4
$ws.Range("AB18").Value = 15.5
$ws.Range("AC18").Value = 65
$ws.Range("AD18").Value = 450
$ws.Range("AF18").Value = 32
$ws.Range("AG18").Value = 17
$ws.Range("AH18").Value = 100
$ws.Range("AI18").Value = 55
$ws.Range("AJ18").Value = 50
$ws.Range("G18").Value = 1.39
$ws.Range("H18").Value = 4.25
$ws.Range("I18").Value = 6.4
$ws.Range("T18").Value = 6
$ws.Range("U18").Value = 5.7
$ws.Range("V18").Value = 7.1
$ws.Range("W18").Value = 7.7
$ws.Range("Y18").Value = 21
# Row 19
$ws.Range("AA19").Value = 15
$ws.Range("AB19").Value = 32
$ws.Range("AC19").Value = 150
$ws.Range("AE19").Value = 40
$ws.Range("AG19").Value = 50
$ws.Range("AH19").Value = 600
$ws.Range("AI19").Value = 200
$ws.Range("AJ19").Value = 150
$ws.Range("G19").Value = 1.1
$ws.Range("H19").Value = 7.2
$ws.Range("I19").Value = 16
$ws.Range("T19").Value = 8.25
$ws.Range("V19").Value = 10
$ws.Range("W19").Value = 5.3
$ws.Range("X19").Value = 9.5
$ws.Range("Y19").Value = 32
# Row 20
$ws.Range("L20").Value = 1.3
$ws.Range("M20").Value = 3.4
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 1.8
# Row 21
$ws.Range("N21").Value = 2.35
$ws.Range("O21").Value = 1.57
# Row 24
$ws.Range("AA24").Value = 13
$ws.Range("AI24").Value = 51
$ws.Range("H24").Value = 6.25
$ws.Range("L24").Value = 1.11
$ws.Range("M24").Value = 6
$ws.Range("R24").Value = 1.83
$ws.Range("S24").Value = 1.83
$ws.Range("T24").Value = 10
$ws.Range("V24").Value = 10
$ws.Range("Z24").Value = 21
# Row 25
$ws.Range("AF25").Value = 21
$ws.Range("AI25").Value = 29
$ws.Range("G25").Value = 1.9
$ws.Range("I25").Value = 3.5
$ws.Range("J25").Value = 1.02
$ws.Range("K25").Value = 12
$ws.Range("L25").Value = 1.22
$ws.Range("M25").Value = 4
$ws.Range("N25").Value = 1.75
$ws.Range("O25").Value = 2.05
$ws.Range("R25").Value = 1.67
$ws.Range("S25").Value = 2.1
# Row 28
$ws.Range("AB28").Value = 16
$ws.Range("AC28").Value = 65
$ws.Range("H28").Value = 4.05
$ws.Range("I28").Value = 5.6
$ws.Range("M28").Value = 4
$ws.Range("Q28").Value = 3.05
$ws.Range("T28").Value = 7.3
$ws.Range("W28").Value = 12.5
$ws.Range("Y28").Value = 24
# Row 29
$ws.Range("AB29").Value = 14.5
$ws.Range("AE29").Value = 10.75
$ws.Range("AG29").Value = 13.5
$ws.Range("AI29").Value = 37
$ws.Range("I29").Value = 3.8
$ws.Range("L29").Value = 1.27
$ws.Range("O29").Value = 1.93
$ws.Range("S29").Value = 2.05
$ws.Range("T29").Value = 7.7
$ws.Range("U29").Value = 10.5
$ws.Range("W29").Value = 18.5
$ws.Range("X29").Value = 15.5
$ws.Range("Y29").Value = 26
# Row 32
$ws.Range("AB32").Value = 26
$ws.Range("AE32").Value = 40
$ws.Range("AF32").Value = 110
$ws.Range("AH32").Value = 400
$ws.Range("H32").Value = 6.8
$ws.Range("I32").Value = 12.5
$ws.Range("N32").Value = 1.32
$ws.Range("O32").Value = 3.1
$ws.Range("Q32").Value = 4.1
$ws.Range("S32").Value = 1.82
$ws.Range("T32").Value = 11.25
$ws.Range("U32").Value = 7.7
$ws.Range("W32").Value = 7.6
$ws.Range("X32").Value = 10.25
$ws.Range("Y32").Value = 26
